$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

Set-TextValue "D2" "24.887.14"
Set-TextValue "E2" "  -4.14%  "
Set-TextValue "D3" "1.632.33"
Set-TextValue "E3" "  -6.39%  "
Set-TextValue "D4" "0.9976"
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "234.64"
Set-TextValue "E5" "  -5.53%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "D7" "0.4725"
Set-TextValue "E7" "  -6.51%  "
Set-TextValue "D8" "0.2547"
Set-TextValue "E8" "  -7.13%  "
Set-TextValue "D9" "0.06089"
Set-TextValue "E9" "  -1.55%  "
Set-TextValue "D10" "0.06934"
Set-TextValue "E10" "  -4.45%  "
Set-TextValue "D11" "1.632.03"
Set-TextValue "E11" "  -6.40%  "
Set-TextValue "D12" "14.63"
Set-TextValue "E12" "  -3.62%  "
Set-TextValue "D13" "0.6101"
Set-TextValue "E13" "  -6.73%  "
Set-TextValue "D14" "4.335"
Set-TextValue "E14" "  -6.62%  "
Set-TextValue "D15" "72.71"
Set-TextValue "E15" "  -6.42%  "
Set-TextValue "E16" "  +0.12%  "
Set-TextValue "D17" "0.9981"
Set-TextValue "E17" "  -0.07%  "
Set-TextValue "D18" "24.902.74"
Set-TextValue "D19" "0.000006554"
Set-TextValue "E19" "  -4.14%  "
Set-TextValue "D20" "11.05"
Set-TextValue "E20" "  -6.60%  "
Set-TextValue "D21" "1.842.35"
Set-TextValue "E21" "  -6.55%  "
Set-TextValue "D22" "4.335"
Set-TextValue "E22" "  -1.08%  "
Set-TextValue "D23" "8.537"
Set-TextValue "E23" "  -1.80%  "
Set-TextValue "D24" "5.233"
Set-TextValue "E24" "  -3.13%  "
Set-TextValue "D25" "133.94"
Set-TextValue "E25" "  -2.00%  "
Set-TextValue "D26" "14.71"
Set-TextValue "E26" "  -3.50%  "
Set-TextValue "E27" "  -8.70%  "
Set-TextValue "D28" "102.56"
Set-TextValue "E28" "  -2.76%  "
Set-TextValue "E29" "  -8.30%  "
Set-TextValue "D30" "3.741"
Set-TextValue "E30" "  -4.01%  "
Set-TextValue "D31" "0.07707"
Set-TextValue "E31" "  -6.46%  "
Set-TextValue "D32" "3.528"
Set-TextValue "E32" "  -2.99%  "
Set-TextValue "D33" "0.9987"
Set-TextValue "E33" "  +0.05%  "
Set-TextValue "D34" "0.04276"
Set-TextValue "E34" "  -8.64%  "
Set-TextValue "D35" "2.596"
Set-TextValue "E35" "  -2.23%  "
Set-TextValue "D36" "0.9167"
Set-TextValue "E36" "  -7.69%  "
Set-TextValue "E37" "  -6.84%  "
Set-TextValue "D38" "2.550"
Set-TextValue "E38" "  -7.35%  "
Set-TextValue "D39" "0.01540"
Set-TextValue "E39" "  -4.46%  "
Set-TextValue "E40" "  -0.12%  "
Set-TextValue "D41" "0.8177"
Set-TextValue "E41" "  +7.76%  "
Set-TextValue "D42" "97.14"
Set-TextValue "E42" "  -2.91%  "
Set-TextValue "D43" "1.768"
Set-TextValue "E43" "  -8.25%  "
Set-TextValue "D44" "0.3674"
Set-TextValue "E44" "  -6.12%  "
Set-TextValue "D45" "4.697"
Set-TextValue "E45" "  -6.26%  "
Set-TextValue "B46" "Cronos"
Set-TextValue "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.05196"
Set-TextValue "E46" "  -1.24%  "
Set-TextValue "B47" "Algorand"
Set-TextValue "C47" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D47" "0.1087"
Set-TextValue "E47" "  -5.09%  "
Set-TextValue "D48" "5.994"
Set-TextValue "E48" "  -4.77%  "
Set-TextValue "D49" "29.40"
Set-TextValue "E49" "  -3.97%  "
Set-TextValue "D50" "0.9997"
Set-TextValue "E50" "  -0.17%  "
Set-TextValue "D51" "0.9999"
Set-TextValue "E51" "  -0.16%  "
